$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(20250527, "A1", 16.18, 160.33699999999999),
    @(20250527, "A2", 16.477, 142.797),
    @(20250527, "A3", 15.737, 139.80699999999999),
    @(20250527, "A4", 14.180999999999999, 113.994),
    @(20250527, "A5", 18.452999999999999, 223.53299999999999),
    @(20250527, "B1", 17.169, 181.102),
    @(20250527, "B2", 15.659000000000001, 145.273),
    @(20250527, "B3", 17.751999999999999, 131.476),
    @(20250527, "B4", 16.82, 161.20699999999999),
    @(20250527, "B5", 21.943999999999999, 212.886),
    @(20250527, "C1", 20.908999999999999, 210.125),
    @(20250527, "C2", 14.603, 144.41200000000001),
    @(20250527, "C3", 16.466000000000001, 163.6),
    @(20250527, "C4", 18.945, 218.34800000000001),
    @(20250527, "C5", 17.736000000000001, 176.03),
    @(20250527, "D1", 15.512, 147.923),
    @(20250527, "D2", 15.987, 169.637),
    @(20250527, "D3", 14.522, 134.59800000000001),
    @(20250527, "D4", 18.655999999999999, 216.816),
    @(20250527, "D5", 18.363, 191.185),
    @(20250527, "E1", 15.992000000000001, 128.91),
    @(20250527, "E2", 20.524000000000001, 209.49600000000001),
    @(20250527, "E3", 20.95, 249.12200000000001),
    @(20250527, "E4", 21.321999999999999, 237.16399999999999),
    @(20250527, "E5", 19.579000000000001, 191.27600000000001)
)

$startRow = 77
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$ws.Range("A77:D101").Font.Color = $ws.Range("B2").Font.Color

$ws.Range("C80").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 54 | Out-Null

